# Applies the weekly Fruta/Hortaliza update to rows 4, 5, 6 and 8
# (row 7 is untouched by this commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---------------------------------------------------------------
$ws.Range("D4").Value = 44169
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21000
$ws.Range("S4").Value = 1167

# --- Row 5 ---------------------------------------------------------------
$ws.Range("D5").Value = 44533
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("Q5").Value = "$/caja 10 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1450
$ws.Range("T5").Value = 10

# --- Row 6 ---------------------------------------------------------------
$ws.Range("D6").Value = 44160
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1361

# --- Row 8 ---------------------------------------------------------------
$ws.Range("D8").Value = 44524
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 27000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 27500
$ws.Range("Q8").Value = "$/bandeja 18 kilos"
$ws.Range("R8").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S8").Value = 1528
$ws.Range("T8").Value = 18
